$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.693.67"
$ws.Range("D2").Style = $ws.Range("A1").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.420.96"
$ws.Range("D3").Style = $ws.Range("A1").Style
$ws.Range("E3").Value = "  -1.85%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.49"
$ws.Range("D5").Style = $ws.Range("A1").Style
$ws.Range("E5").Value = "  -1.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.53"
$ws.Range("D6").Style = $ws.Range("A1").Style
$ws.Range("E6").Value = "  -1.96%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.610"
$ws.Range("D7").Style = $ws.Range("A1").Style
$ws.Range("E7").Value = "  +4.39%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.427.00"
$ws.Range("D9").Style = $ws.Range("A1").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.14"
$ws.Range("D10").Style = $ws.Range("A1").Style
$ws.Range("E10").Value = "  -2.44%  "
$ws.Range("E11").Value = "  -3.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.441"
$ws.Range("D12").Style = $ws.Range("A1").Style
$ws.Range("E12").Value = "  -1.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.010.01"
$ws.Range("D13").Style = $ws.Range("A1").Style
$ws.Range("E13").Value = "  -1.91%  "
$ws.Range("E14").Value = "  -0.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000188"
$ws.Range("D15").Style = $ws.Range("A1").Style
$ws.Range("E15").Value = "  -4.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.60"
$ws.Range("D16").Style = $ws.Range("A1").Style
$ws.Range("E16").Value = "  -3.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.689.72"
$ws.Range("D17").Style = $ws.Range("A1").Style
$ws.Range("E17").Value = "  -1.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.411.07"
$ws.Range("D18").Style = $ws.Range("A1").Style
$ws.Range("E18").Value = "  -1.81%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.35"
$ws.Range("D19").Style = $ws.Range("A1").Style
$ws.Range("E19").Value = "  -1.80%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.88"
$ws.Range("D20").Style = $ws.Range("A1").Style
$ws.Range("E20").Value = "  -3.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "379.80"
$ws.Range("D21").Style = $ws.Range("A1").Style
$ws.Range("E21").Value = "  -3.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.00"
$ws.Range("D22").Style = $ws.Range("A1").Style
$ws.Range("E22").Value = "  -3.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.549"
$ws.Range("D23").Style = $ws.Range("A1").Style
$ws.Range("E23").Value = "  -0.77%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.10"
$ws.Range("D25").Style = $ws.Range("A1").Style
$ws.Range("E25").Value = "  -2.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000118"
$ws.Range("D26").Style = $ws.Range("A1").Style
$ws.Range("E26").Value = "  -5.90%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.13"
$ws.Range("D27").Style = $ws.Range("A1").Style
$ws.Range("E27").Value = "  +5.76%  "
$ws.Range("E28").Value = "  -1.40%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = $ws.Range("A1").Style
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.47"
$ws.Range("D30").Style = $ws.Range("A1").Style
$ws.Range("E30").Value = "  +1.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.15"
$ws.Range("D31").Style = $ws.Range("A1").Style
$ws.Range("E31").Value = "  -4.79%  "
$ws.Range("E32").Value = "  -2.44%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.25"
$ws.Range("D33").Style = $ws.Range("A1").Style
$ws.Range("E33").Value = "  -2.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.07"
$ws.Range("D34").Style = $ws.Range("A1").Style
$ws.Range("E34").Value = "  -1.31%  "
$ws.Range("E35").Value = "  +3.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "160.59"
$ws.Range("D36").Style = $ws.Range("A1").Style
$ws.Range("E36").Value = "  -0.46%  "
$ws.Range("E37").Value = "  -3.14%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0758"
$ws.Range("D38").Style = $ws.Range("A1").Style
$ws.Range("E38").Value = "  -2.37%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.84"
$ws.Range("D39").Style = $ws.Range("A1").Style
$ws.Range("E39").Value = "  +2.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.888.71"
$ws.Range("D40").Style = $ws.Range("A1").Style
$ws.Range("E40").Value = "  -5.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.46"
$ws.Range("D41").Style = $ws.Range("A1").Style
$ws.Range("E41").Value = "  -3.22%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.60"
$ws.Range("D42").Style = $ws.Range("A1").Style
$ws.Range("E42").Value = "  +1.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "43.01"
$ws.Range("D43").Style = $ws.Range("A1").Style
$ws.Range("E43").Value = "  +0.42%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0315"
$ws.Range("D44").Style = $ws.Range("A1").Style
$ws.Range("E44").Value = "  -2.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.772"
$ws.Range("D45").Style = $ws.Range("A1").Style
$ws.Range("E45").Value = "  -0.75%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.79"
$ws.Range("D46").Style = $ws.Range("A1").Style
$ws.Range("E46").Value = "  -0.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "320.07"
$ws.Range("D47").Style = $ws.Range("A1").Style
$ws.Range("E47").Value = "  +2.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.06"
$ws.Range("D48").Style = $ws.Range("A1").Style
$ws.Range("E48").Value = "  -6.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.20"
$ws.Range("D49").Style = $ws.Range("A1").Style
$ws.Range("E49").Value = "  -2.73%  "
$ws.Range("E50").Value = "  +0.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.55"
$ws.Range("D51").Style = $ws.Range("A1").Style
$ws.Range("E51").Value = "  -2.88%  "
